# [Fonds de solidarite] Add 2021-03-17 data
#
# Updates a handful of "nombre_aides" / "nombre_entreprises" / "montant_total"
# figures (all stored as text in the sheet) for:
#   - Auvergne-Rhone-Alpes / SARL   (row 4)
#   - Auvergne-Rhone-Alpes / SAS    (row 6)
#   - Corse / SARL                  (row 28)
#   - Occitanie / SARL              (row 88)
#   - Occitanie / SAS               (row 90)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell address -> new (text) value. All of these cells already hold plain
# numeric-looking text (e.g. "1608", "13979280.60"), so we must force the
# cell to stay text (otherwise Excel auto-converts it to a real number and
# mangles values like the trailing zero in "13979280.60").
$changes = [ordered]@{
    "C4"  = "1610"          # nombre_aides
    "D4"  = "1382"          # nombre_entreprises
    "E4"  = "14007320.94"   # montant_total

    "C6"  = "1120"          # nombre_aides
    "D6"  = "977"           # nombre_entreprises
    "E6"  = "8582679.67"    # montant_total

    "C28" = "182"           # nombre_aides
    "E28" = "1151374.44"    # montant_total

    "C88" = "1534"          # nombre_aides
    "E88" = "11778274.14"   # montant_total

    "C90" = "1098"          # nombre_aides
    "D90" = "941"           # nombre_entreprises
    "E90" = "7799142.44"    # montant_total
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    # Force a text number format so the assigned string isn't reinterpreted
    # as a numeric value (which would both change the cell type and risk
    # floating point rounding of the decimal amounts).
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
    # Drop the temporary text format again so the cell's style matches the
    # original (un-styled / default) cells as closely as possible.
    $cell.ClearFormats()
}
